# Auto-generated Excel COM-interop script applying the Gilgamesh_Profits
# leve-profit recalculation update (currentAveragePrice* / Leve*Price* /
# LeveProfit* columns) across all eight crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1097.875
$ws.Range("I17").Value = 1599.8
$ws.Range("J17").Value = 965.7895
$ws.Range("K17").Value = 4799.4
$ws.Range("L17").Value = 2897.3685
$ws.Range("M17").Value = -4631.4
$ws.Range("N17").Value = -3233.3685
$ws.Range("H41").Value = 445.875
$ws.Range("I41").Value = 392.25
$ws.Range("J41").Value = 499.5
$ws.Range("K41").Value = 392.25
$ws.Range("L41").Value = 499.5
$ws.Range("M41").Value = 47.75
$ws.Range("N41").Value = -1379.5
$ws.Range("H64").Value = 250004670
$ws.Range("I64").Value = 4999.5
$ws.Range("K64").Value = 4999.5
$ws.Range("M64").Value = -4751.5
$ws.Range("H67").Value = 250004670
$ws.Range("I67").Value = 4999.5
$ws.Range("K67").Value = 4999.5
$ws.Range("M67").Value = -4141.5
$ws.Range("H70").Value = 1668
$ws.Range("J70").Value = 1668
$ws.Range("L70").Value = 5004
$ws.Range("N70").Value = -5544
$ws.Range("H73").Value = 1668
$ws.Range("J73").Value = 1668
$ws.Range("L73").Value = 5004
$ws.Range("N73").Value = -6876
$ws.Range("H111").Value = 2910.7222
$ws.Range("J111").Value = 1422.375
$ws.Range("L111").Value = 4267.125
$ws.Range("N111").Value = -10401.125
$ws.Range("H135").Value = 353.85715
$ws.Range("I135").Value = 372.6154
$ws.Range("J135").Value = 110
$ws.Range("K135").Value = 3353.5386
$ws.Range("L135").Value = 990
$ws.Range("M135").Value = -818.5386000000003
$ws.Range("N135").Value = -6060
$ws.Range("H137").Value = 1828.909
$ws.Range("I137").Value = 1653.375
$ws.Range("K137").Value = 4960.125
$ws.Range("M137").Value = -2410.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1599.4286
$ws.Range("I2").Value = 1833
$ws.Range("K2").Value = 1833
$ws.Range("M2").Value = -1720
$ws.Range("H61").Value = 2774.8
$ws.Range("I61").Value = 1749.9
$ws.Range("K61").Value = 1749.9
$ws.Range("M61").Value = -1537.9
$ws.Range("H74").Value = 2026.12
$ws.Range("I74").Value = 1463
$ws.Range("J74").Value = 2870.8
$ws.Range("K74").Value = 1463
$ws.Range("L74").Value = 2870.8
$ws.Range("M74").Value = -589
$ws.Range("N74").Value = -4618.8
$ws.Range("H77").Value = 2026.12
$ws.Range("I77").Value = 1463
$ws.Range("J77").Value = 2870.8
$ws.Range("K77").Value = 7315
$ws.Range("L77").Value = 14354
$ws.Range("M77").Value = -2947
$ws.Range("N77").Value = -23090
$ws.Range("H116").Value = 1599.4286
$ws.Range("I116").Value = 1833
$ws.Range("K116").Value = 1833
$ws.Range("M116").Value = 461
$ws.Range("H132").Value = 2723.8096
$ws.Range("J132").Value = 2990.75
$ws.Range("L132").Value = 8972.25
$ws.Range("N132").Value = -14032.25
$ws.Range("H135").Value = 106995.5
$ws.Range("J135").Value = 106995.5
$ws.Range("L135").Value = 106995.5
$ws.Range("N135").Value = -117135.5
$ws.Range("H136").Value = 2774.8
$ws.Range("I136").Value = 1749.9
$ws.Range("K136").Value = 5249.700000000001
$ws.Range("M136").Value = -2699.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1599.4286
$ws.Range("I3").Value = 1833
$ws.Range("K3").Value = 1833
$ws.Range("M3").Value = -1719
$ws.Range("H107").Value = 2025533.5
$ws.Range("I107").Value = 2653693.2
$ws.Range("J107").Value = 1462.7778
$ws.Range("K107").Value = 2653693.2
$ws.Range("L107").Value = 1462.7778
$ws.Range("M107").Value = -2651773.2
$ws.Range("N107").Value = -5302.7778
$ws.Range("H132").Value = 109998
$ws.Range("J132").Value = 109998
$ws.Range("L132").Value = 109998
$ws.Range("N132").Value = -120118
$ws.Range("H134").Value = 2087.5454
$ws.Range("J134").Value = 3182.889
$ws.Range("L134").Value = 9548.667000000001
$ws.Range("N134").Value = -14618.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5558666.5
$ws.Range("I6").Value = 16670000
$ws.Range("K6").Value = 16670000
$ws.Range("M6").Value = -16669887
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H62").Value = 3999.2
$ws.Range("I62").Value = 3998.5
$ws.Range("K62").Value = 3998.5
$ws.Range("M62").Value = -3374.5
$ws.Range("H65").Value = 3999.2
$ws.Range("I65").Value = 3998.5
$ws.Range("K65").Value = 19992.5
$ws.Range("M65").Value = -16872.5
$ws.Range("H132").Value = 3807.2354
$ws.Range("I132").Value = 3182.5386
$ws.Range("K132").Value = 9547.6158
$ws.Range("M132").Value = -7017.6158
$ws.Range("H134").Value = 3745.4119
$ws.Range("I134").Value = 4259.92
$ws.Range("K134").Value = 12779.76
$ws.Range("M134").Value = -10244.76

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 17883.25
$ws.Range("I56").Value = 17883.25
$ws.Range("K56").Value = 17883.25
$ws.Range("M56").Value = -17353.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 11000
$ws.Range("J24").Value = 11000
$ws.Range("L24").Value = 11000
$ws.Range("N24").Value = -11346
$ws.Range("H80").Value = 111116080
$ws.Range("I80").Value = 333336670
$ws.Range("J80").Value = 5785.3335
$ws.Range("K80").Value = 333336670
$ws.Range("L80").Value = 5785.3335
$ws.Range("M80").Value = -333335672
$ws.Range("N80").Value = -7781.3335
$ws.Range("H83").Value = 111116080
$ws.Range("I83").Value = 333336670
$ws.Range("J83").Value = 5785.3335
$ws.Range("K83").Value = 1666683350
$ws.Range("L83").Value = 28926.6675
$ws.Range("M83").Value = -1666678358
$ws.Range("N83").Value = -38910.6675
$ws.Range("H97").Value = 1624.3572
$ws.Range("J97").Value = 974.5
$ws.Range("L97").Value = 974.5
$ws.Range("N97").Value = -1966.5
$ws.Range("H102").Value = 1250.3024
$ws.Range("I102").Value = 913.35
$ws.Range("J102").Value = 1543.3043
$ws.Range("K102").Value = 913.35
$ws.Range("L102").Value = 1543.3043
$ws.Range("M102").Value = 708.65
$ws.Range("N102").Value = -4787.3043
$ws.Range("H126").Value = 7851.5415
$ws.Range("I126").Value = 2050.5557
$ws.Range("J126").Value = 11332.134
$ws.Range("K126").Value = 6151.6671
$ws.Range("L126").Value = 33996.402
$ws.Range("M126").Value = -3681.6671
$ws.Range("N126").Value = -38936.402
$ws.Range("H132").Value = 2474.923
$ws.Range("I132").Value = 2117
$ws.Range("J132").Value = 2540
$ws.Range("K132").Value = 6351
$ws.Range("L132").Value = 7620
$ws.Range("M132").Value = -3821
$ws.Range("N132").Value = -12680
$ws.Range("H136").Value = 73663
$ws.Range("J136").Value = 73663
$ws.Range("L136").Value = 220989
$ws.Range("N136").Value = -226089

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1979.1111
$ws.Range("I7").Value = 1393.5454
$ws.Range("J7").Value = 2899.2856
$ws.Range("K7").Value = 1393.5454
$ws.Range("L7").Value = 2899.2856
$ws.Range("M7").Value = -1281.5454
$ws.Range("N7").Value = -3123.2856
$ws.Range("H16").Value = 1959.6666
$ws.Range("I16").Value = 1944.5
$ws.Range("K16").Value = 1944.5
$ws.Range("M16").Value = -1774.5
$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2705
$ws.Range("H27").Value = 3000
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -2893
$ws.Range("H107").Value = 4997
$ws.Range("I107").Value = 4997
$ws.Range("K107").Value = 4997
$ws.Range("M107").Value = -3077
$ws.Range("H126").Value = 1979.1111
$ws.Range("I126").Value = 1393.5454
$ws.Range("J126").Value = 2899.2856
$ws.Range("K126").Value = 4180.6362
$ws.Range("L126").Value = 8697.856800000001
$ws.Range("M126").Value = -1710.6362
$ws.Range("N126").Value = -13637.8568
$ws.Range("H141").Value = 80000
$ws.Range("I141").Value = 80000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 80000
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0
$ws.Range("M141").Value = -74820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1246.1818
$ws.Range("I2").Value = 1270.8
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1270.8
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -1158.8
$ws.Range("N2").Value = -1224
$ws.Range("H4").Value = 302.1
$ws.Range("I4").Value = 257
$ws.Range("J4").Value = 407.33334
$ws.Range("K4").Value = 257
$ws.Range("L4").Value = 407.33334
$ws.Range("M4").Value = -144
$ws.Range("N4").Value = -633.33334
$ws.Range("H49").Value = 34999
$ws.Range("I49").Value = 34999
$ws.Range("K49").Value = 34999
$ws.Range("M49").Value = -34769
$ws.Range("H132").Value = 3268.72
$ws.Range("I132").Value = 3266.1365
$ws.Range("J132").Value = 3287.6667
$ws.Range("K132").Value = 9798.4095
$ws.Range("L132").Value = 9863.000100000001
$ws.Range("M132").Value = -7268.4095
$ws.Range("N132").Value = -14923.0001
$ws.Range("H136").Value = 3971.3572
$ws.Range("I136").Value = 4300.2
$ws.Range("J136").Value = 3149.25
$ws.Range("K136").Value = 12900.6
$ws.Range("L136").Value = 9447.75
$ws.Range("M136").Value = -10350.6
$ws.Range("N136").Value = -14547.75

